# Reorder slides: the slide "How to Run the Game" (currently at position 10)
# and the slide "How the Game is Scored" (currently at position 12) swap
# places, while the slide in between ("Tips to Test Your Code", position 11)
# stays put.
#
# A plain MoveTo(10 -> 12) would shift position 11 up to 10, which is not
# what we want, so we move the slide at position 10 to 12 first (pushing the
# in-between slide back to position 10), then move that in-between slide
# (now sitting at position 11) back down to position 10 - producing a clean
# pairwise swap of positions 10 and 12 only.

$p = $ppt.ActivePresentation

$runGame = $p.Slides.Item(10)     # "How to Run the Game"
$runGame.MoveTo(12)

$tips = $p.Slides.Item(11)        # "Tips to Test Your Code", shifted up to 11
$tips.MoveTo(10)
